# "Generate Report for Archive"
#
# The localization status report is regenerated: the two still-untranslated
# source files flip from "Ready for handoff" to "In Translation" everywhere
# that status is shown (the Overview roll-up columns for each locale, plus
# the per-locale "Status" column on the zh-cn / de-de detail sheets). Excel
# then re-narrows those status columns to fit the new, shorter label.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---------------
# Overview sheet: one status column per locale (E = zh-cn, F = de-de)
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# Per-locale detail sheets: "Status" is column C
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Column widths follow the new (shorter) status text ------------------
$newStatusWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newStatusWidth   # E:E
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth   # F:F
$zhcn.Columns.Item(3).ColumnWidth     = $newStatusWidth   # C:C
$dede.Columns.Item(3).ColumnWidth     = $newStatusWidth   # C:C
